# Updates the "想去人数" (F column) values on the "展览" and "全部类型"
# worksheets to match the newly generated gh-pages data dump.

$wb = $excel.ActiveWorkbook

# row -> new F-column value, identical for both "展览" and "全部类型" sheets
# except for row 45, which already differed by one between the two sheets
# before the edit (1886 vs 1887) and converges to the same value (1994).
$commonUpdates = @{
    3  = 1467
    4  = 180
    6  = 244
    7  = 107
    9  = 203
    10 = 154
    12 = 4872
    14 = 7162
    21 = 4218
    22 = 1613
    24 = 85
    25 = 2799
    28 = 187
    29 = 417
    30 = 404
    31 = 423
    32 = 260
    34 = 1659
    35 = 1103
    37 = 1205
    38 = 96
    39 = 560
    40 = 14
    41 = 505
    43 = 36
    44 = 98
    46 = 668
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonUpdates.Keys) {
        $ws.Range("F$row").Value = $commonUpdates[$row]
    }

    # row 45 goes to 1994 on both sheets
    $ws.Range("F45").Value = 1994
}
